# W13 Wednesday Commit 1
# Fills in the week-2 (rows 6-11) observations for all four sheets, adds a
# new "Leaf 6" column (with data) to the two Cardoon sheets, highlights a
# couple of outlier cells, restores each sheet's selection, and sets
# portrait orientation on the last sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Pansies Alive  (B6:F11)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Pansies Alive")
$pansiesAlive = @(
    @(0, 1, 2, 2, 3),
    @(0, 0, 2, 2, 3),
    @(1, 1, 3, 2, 3),
    @(1, 1, 3, 2, 3),
    @(1, 2, 3, 3, 3),
    @(1, 2, 4, 3, 3)
)
for ($i = 0; $i -lt $pansiesAlive.Count; $i++) {
    $row = 6 + $i
    $vals = $pansiesAlive[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws1.Cells.Item($row, 2 + $j).Value = $vals[$j]
    }
}

# ---------------------------------------------------------------------
# Sheet: Pansies Dead  (B6:F11)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Pansies Dead")
$pansiesDead = @(
    @(0, 0, 1, 0, 1),
    @(0, 1, 2, 0, 0),
    @(0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 1)
)
for ($i = 0; $i -lt $pansiesDead.Count; $i++) {
    $row = 6 + $i
    $vals = $pansiesDead[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws2.Cells.Item($row, 2 + $j).Value = $vals[$j]
    }
}

# ---------------------------------------------------------------------
# Sheet: Cardoon (2)  -- create the red-font highlight FIRST so the new
# font lands at fontId=1 (red) ahead of the amber font at fontId=2, which
# matches the workbook's recorded style order.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Cardoon (2)")
$ws4.Range("C11").Value = 13
$ws4.Range("C11").Font.Color = 255

# ---------------------------------------------------------------------
# Sheet: Cardoon (1)  (B6:F11, new column G header + G9:G11)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Cardoon (1)")
$ws3.Range("G1").Value = "Leaf 6"

$cardoon1 = @(
    @(19,   16,    8,    22,    12.5),
    @(19.5, 17,    8,    22,    12.5),
    @(20,   17,    8,    22.5,  12.5),
    @(19.5, 18,    7.5,  22.5,  13),
    @(20,   18.5,  7.5,  21.75, 12.5),
    @(20,   19.34, 7.25, 21.5,  12.5)
)
for ($i = 0; $i -lt $cardoon1.Count; $i++) {
    $row = 6 + $i
    $vals = $cardoon1[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws3.Cells.Item($row, 2 + $j).Value = $vals[$j]
    }
}
# D9 is the amber-highlighted outlier -- creates fontId=2 (amber)
$ws3.Range("D9").Value = 7.5
$ws3.Range("D9").Font.Color = 49407

$ws3.Range("G9").Value = 5
$ws3.Range("G10").Value = 6
$ws3.Range("G11").Value = 6.75

# ---------------------------------------------------------------------
# Sheet: Cardoon (2)  (B6:F11, new column G header + G9:G11)
# ---------------------------------------------------------------------
$ws4.Range("G1").Value = "Leaf 6"

$cardoon2 = @(
    @(24,    13.5,  18.5,  24.5, 24),
    @(24,    13,    20.75, 24,   24),
    @(24,    13.5,  21.5,  23.75,24.5),
    @(24.5,  13,    22.5,  24,   24.5),
    @(24.75, 13.25, 23.25, 24,   24.5),
    @(24.25, 13,    24.18, 24,   24.5)
)
for ($i = 0; $i -lt $cardoon2.Count; $i++) {
    $row = 6 + $i
    $vals = $cardoon2[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws4.Cells.Item($row, 2 + $j).Value = $vals[$j]
    }
}
# C9 is the amber-highlighted outlier (reuses fontId=2 created above)
$ws4.Range("C9").Value = 13
$ws4.Range("C9").Font.Color = 49407

$ws4.Range("G9").Value = 5
$ws4.Range("G10").Value = 6
$ws4.Range("G11").Value = 7

# Page orientation for the last (active) sheet
$ws4.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Restore each sheet's selection without leaving a different tab active.
# ---------------------------------------------------------------------
$excel.Goto($ws1.Range("G18"))
$excel.Goto($ws2.Range("N20"))
$excel.Goto($ws3.Range("J9"))
$excel.Goto($ws4.Range("F17"))
$ws4.Activate()
